$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2620.7144
$ws.Range("J19").Value = 3431.6
$ws.Range("L19").Value = 3431.6
$ws.Range("N19").Value = -3781.6
$ws.Range("H39").Value = 273.3871
$ws.Range("I39").Value = 168
$ws.Range("K39").Value = 504
$ws.Range("M39").Value = -208
$ws.Range("H80").Value = 1122.5
$ws.Range("I80").Value = 646.75
$ws.Range("K80").Value = 1940.25
$ws.Range("M80").Value = -942.25
$ws.Range("H83").Value = 1122.5
$ws.Range("I83").Value = 646.75
$ws.Range("K83").Value = 5820.75
$ws.Range("M83").Value = -828.75
$ws.Range("H127").Value = 4057.2
$ws.Range("I127").Value = 662
$ws.Range("K127").Value = 1986
$ws.Range("M127").Value = 2974

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 50000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H39").Value = 5008
$ws.Range("I39").Value = 5008
$ws.Range("K39").Value = 5008
$ws.Range("M39").Value = -4488
$ws.Range("H40").Value = 21000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 6363.1816
$ws.Range("I41").Value = 5555
$ws.Range("K41").Value = 5555
$ws.Range("M41").Value = -5141
$ws.Range("H44").Value = 42800.6
$ws.Range("J44").Value = 42800.6
$ws.Range("L44").Value = 42800.6
$ws.Range("N44").Value = -43776.6
$ws.Range("H74").Value = 7651031.5
$ws.Range("I74").Value = 11366050
$ws.Range("K74").Value = 11366050
$ws.Range("M74").Value = -11365176
$ws.Range("H77").Value = 7651031.5
$ws.Range("I77").Value = 11366050
$ws.Range("K77").Value = 56830250
$ws.Range("M77").Value = -56825882
$ws.Range("H132").Value = 2783.012
$ws.Range("I132").Value = 1995.164
$ws.Range("K132").Value = 5985.492
$ws.Range("M132").Value = -3455.492

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 700
$ws.Range("I19").Value = 700
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 700
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -527
$ws.Range("H25").Value = 424.75
$ws.Range("J25").Value = 650
$ws.Range("L25").Value = 650
$ws.Range("N25").Value = -1120
$ws.Range("H46").Value = 6999
$ws.Range("J46").Value = 6999
$ws.Range("L46").Value = 6999
$ws.Range("N46").Value = -7595

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4000
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H50").Value = 44999.832
$ws.Range("J50").Value = 44999.832
$ws.Range("L50").Value = 44999.832
$ws.Range("N50").Value = -46249.832
$ws.Range("H60").Value = 6500
$ws.Range("I60").Value = 6500
$ws.Range("K60").Value = 6500
$ws.Range("M60").Value = -5989
$ws.Range("H62").Value = 1253927.5
$ws.Range("J62").Value = 2899.4
$ws.Range("L62").Value = 2899.4
$ws.Range("N62").Value = -4147.4
$ws.Range("H65").Value = 1253927.5
$ws.Range("J65").Value = 2899.4
$ws.Range("L65").Value = 14497
$ws.Range("N65").Value = -20737

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 228.67741
$ws.Range("J15").Value = 234.79167
$ws.Range("L15").Value = 704.37501
$ws.Range("N15").Value = -984.37501
$ws.Range("H19").Value = 145
$ws.Range("J19").Value = 100
$ws.Range("L19").Value = 300
$ws.Range("N19").Value = -648
$ws.Range("H21").Value = 459.2
$ws.Range("J21").Value = 994
$ws.Range("L21").Value = 2982
$ws.Range("N21").Value = -3328
$ws.Range("H26").Value = 189.83333
$ws.Range("J26").Value = 153
$ws.Range("L26").Value = 459
$ws.Range("N26").Value = -1035
$ws.Range("H33").Value = 115.38461
$ws.Range("I33").Value = 14.166667
$ws.Range("J33").Value = 202.14285
$ws.Range("K33").Value = 85.00000199999999
$ws.Range("L33").Value = 1212.8571
$ws.Range("M33").Value = 197.999998
$ws.Range("N33").Value = -1778.8571
$ws.Range("H46").Value = 1230.8334
$ws.Range("I46").Value = 346.25
$ws.Range("K46").Value = 1038.75
$ws.Range("M46").Value = -947.75
$ws.Range("H49").Value = 9670.333000000001
$ws.Range("I49").Value = 9003
$ws.Range("J49").Value = 10004
$ws.Range("K49").Value = 27009
$ws.Range("L49").Value = 30012
$ws.Range("M49").Value = -26853
$ws.Range("N49").Value = -30324
$ws.Range("H97").Value = 11906296
$ws.Range("J97").Value = 2302
$ws.Range("L97").Value = 6906
$ws.Range("N97").Value = -7898
$ws.Range("H131").Value = 4924.5083
$ws.Range("I131").Value = 4202.7144
$ws.Range("J131").Value = 5149.067
$ws.Range("K131").Value = 12608.1432
$ws.Range("L131").Value = 15447.201
$ws.Range("M131").Value = -7568.143199999999
$ws.Range("N131").Value = -25527.201

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 750
$ws.Range("I9").Value = 1250
$ws.Range("K9").Value = 1250
$ws.Range("M9").Value = -1080
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H41").Value = 14999.667
$ws.Range("J41").Value = 14999.667
$ws.Range("L41").Value = 14999.667
$ws.Range("N41").Value = -15709.667
$ws.Range("H45").Value = 45916.668
$ws.Range("J45").Value = 45916.668
$ws.Range("L45").Value = 45916.668
$ws.Range("N45").Value = -47034.668

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4626.407
$ws.Range("I46").Value = 2015.65
$ws.Range("K46").Value = 2015.65
$ws.Range("M46").Value = -1827.65
$ws.Range("H105").Value = 120000
$ws.Range("J105").Value = 120000
$ws.Range("L105").Value = 120000
$ws.Range("N105").Value = -126988

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H39").Value = 30490
$ws.Range("I39").Value = 30485
$ws.Range("K39").Value = 30485
$ws.Range("M39").Value = -30072
$ws.Range("H42").Value = 59999
$ws.Range("I42").Value = 59999
$ws.Range("K42").Value = 59999
$ws.Range("M42").Value = -59621
$ws.Range("H43").Value = 77450
$ws.Range("I43").Value = 74900
$ws.Range("K43").Value = 74900
$ws.Range("M43").Value = -74751

